$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): updated coin prices ---

# Some new prices (e.g. "0.999", "4.92") would otherwise be auto-detected as
# numbers by Excel, so force the Text number format before writing those so they
# stay plain text, matching the rest of the sheet.
$priceTextUpdates = [ordered]@{
    'D5' = '583.63'
    'D6' = '171.54'
    'D8' = '0.513'
    'D12' = '4.92'
    'D15' = '25.34'
    'D17' = '0.0000170'
    'D19' = '11.03'
    'D20' = '7.40'
    'D21' = '348.26'
    'D22' = '4.03'
    'D23' = '0.999'
    'D24' = '68.51'
    'D26' = '1.80'
    'D27' = '9.30'
    'D31' = '511.93'
    'D32' = '7.74'
    'D33' = '1.23'
    'D36' = '160.07'
    'D39' = '18.25'
    'D44' = '4.81'
    'D46' = '38.76'
    'D47' = '142.58'
}
foreach ($cellRef in $priceTextUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceTextUpdates[$cellRef]
}

# These new prices (e.g. "67.122.55") already contain multiple separators, so
# Excel keeps them as text without any extra help.
$priceUpdates = [ordered]@{
    'D2' = '67.122.55'
    'D3' = '2.487.75'
    'D9' = '2.487.41'
    'D14' = '2.973.48'
    'D16' = '67.115.53'
    'D18' = '2.460.50'
    'D29' = '2.615.06'
    'D30' = '0.0₃0906'
}
foreach ($cellRef in $priceUpdates.Keys) {
    $ws.Range($cellRef).Value = $priceUpdates[$cellRef]
}

# --- Volume(1h) column (E): updated 1h change percentages ---
# Values are already non-numeric-looking text (percent with surrounding
# spaces), so they can be assigned directly.
$volumeUpdates = [ordered]@{
    'E2' = '  -0.18%  '
    'E3' = '  +0.05%  '
    'E4' = '  +0.01%  '
    'E5' = '  -0.20%  '
    'E6' = '  +2.95%  '
    'E7' = '  -0.06%  '
    'E8' = '  -0.69%  '
    'E9' = '  +0.01%  '
    'E10' = '  +0.71%  '
    'E11' = '  -0.10%  '
    'E12' = '  -0.45%  '
    'E13' = '  -2.17%  '
    'E14' = '  +1.12%  '
    'E15' = '  -2.27%  '
    'E16' = '  +0.17%  '
    'E17' = '  -1.85%  '
    'E18' = '  +1.04%  '
    'E19' = '  -5.24%  '
    'E20' = '  -5.45%  '
    'E21' = '  -3.36%  '
    'E22' = '  -2.16%  '
    'E23' = '  -0.17%  '
    'E24' = '  -2.90%  '
    'E25' = '  -4.25%  '
    'E26' = '  -2.83%  '
    'E27' = '  -1.24%  '
    'E28' = '  +0.21%  '
    'E29' = '  +0.11%  '
    'E30' = '  -3.17%  '
    'E31' = '  +2.59%  '
    'E32' = '  -3.72%  '
    'E33' = '  -2.98%  '
    'E34' = '  -3.79%  '
    'E35' = '  -0.04%  '
    'E36' = '  +0.44%  '
    'E37' = '  -7.51%  '
    'E38' = '  +0.64%  '
    'E39' = '  -4.33%  '
    'E40' = '  -5.85%  '
    'E41' = '  -2.41%  '
    'E42' = '  -0.08%  '
    'E43' = '  -2.17%  '
    'E44' = '  -2.80%  '
    'E45' = '  -4.22%  '
    'E46' = '  -1.47%  '
    'E47' = '  +0.54%  '
    'E48' = '  -4.21%  '
    'E49' = '  -4.80%  '
    'E50' = '  -5.09%  '
    'E51' = '  -0.70%  '
}
foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}
